$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "0.39.0"
$ws.Range("C5").Value = "5.5.0"
$ws.Range("C7").Value = "8.8.9"
$ws.Range("C8").Value = "8.0.0"
$ws.Range("C9").Value = "110.4.458"
$ws.Range("C13").Value = "12.10"
$ws.Range("C15").Value = "2020-11-20 12:39"
$ws.Range("C16").Value = "83.0"
$ws.Range("C17").Value = "2.29.2"
$ws.Range("C20").Value = "6.34"
$ws.Range("C21").Value = "12.11.0.26"
$ws.Range("C22").Value = "15.8.7"
$ws.Range("C23").Value = "51.0.0"
$ws.Range("C24").Value = "26.0.2"
$ws.Range("C25").Value = "4.4.160"
$ws.Range("C26").Value = "3.9.0"
$ws.Range("C27").Value = "3.13"
$ws.Range("C28").Value = "1.146.916.0"
$ws.Range("C30").Value = "1.51"
$ws.Range("C31").Value = "5.17.9"
$ws.Range("C32").Value = "15.0.1"
